{"js": "const replacements = [\n  [\"2024-09-12 Thursday\", \"2024-09-13 Friday\"],\n  [\"448\u00d78=3584\", \"816\u00d74=3264\"],\n  [\"556\u00d79=5004\", \"460\u00d78=3680\"],\n  [\"799\u00d77=5593\", \"346\u00d76=2076\"],\n  [\"584\u00d73=1752\", \"495\u00d78=3960\"],\n  [\"376\u00d73=1128\", \"702\u00d76=4212\"],\n  [\"223\u00d79=2007\", \"749\u00d72=1498\"],\n  [\"281\u00d73=843\", \"614\u00d73=1842\"],\n  [\"945\u00d72=1890\", \"265\u00d77=1855\"],\n  [\"450\u00d74=1800\", \"838\u00d79=7542\"],\n  [\"733\u00d77=5131\", \"916\u00d72=1832\"],\n  [\"496\u00d74=1984\", \"845\u00d79=7605\"],\n  [\"102\u00d75=510\", \"180\u00d76=1080\"],\n  [\"772\u00d72=1544\", \"682\u00d73=2046\"],\n  [\"452\u00d78=3616\", \"583\u00d72=1166\"],\n  [\"164\u00d73=492\", \"758\u00d78=6064\"],\n  [\"557\u00d74=2228\", \"587\u00d75=2935\"],\n  [\"628\u00d73=1884\", \"723\u00d72=1446\"],\n  [\"697\u00d79=6273\", \"611\u00d72=1222\"],\n  [\"433\u00d76=2598\", \"376\u00d72=752\"],\n  [\"207\u00d79=1863\", \"774\u00d72=1548\"],\n  [\"485\u00d73=1455\", \"372\u00d73=1116\"],\n  [\"519\u00d72=1038\", \"382\u00d76=2292\"],\n  [\"982\u00d78=7856\", \"906\u00d77=6342\"],\n  [\"298\u00d73=894\", \"810\u00d74=3240\"],\n  [\"369\u00d76=2214\", \"298\u00d76=1788\"],\n];\n\nconst body = context.document.body;\n\nfor (const [from, to] of replacements) {\n  const found = body.search(from, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-09-12 Thursday\", \"2024-09-13 Friday\"),\n    @(\"448\u00d78=3584\", \"816\u00d74=3264\"),\n    @(\"556\u00d79=5004\", \"460\u00d78=3680\"),\n    @(\"799\u00d77=5593\", \"346\u00d76=2076\"),\n    @(\"584\u00d73=1752\", \"495\u00d78=3960\"),\n    @(\"376\u00d73=1128\", \"702\u00d76=4212\"),\n    @(\"223\u00d79=2007\", \"749\u00d72=1498\"),\n    @(\"281\u00d73=843\", \"614\u00d73=1842\"),\n    @(\"945\u00d72=1890\", \"265\u00d77=1855\"),\n    @(\"450\u00d74=1800\", \"838\u00d79=7542\"),\n    @(\"733\u00d77=5131\", \"916\u00d72=1832\"),\n    @(\"496\u00d74=1984\", \"845\u00d79=7605\"),\n    @(\"102\u00d75=510\", \"180\u00d76=1080\"),\n    @(\"772\u00d72=1544\", \"682\u00d73=2046\"),\n    @(\"452\u00d78=3616\", \"583\u00d72=1166\"),\n    @(\"164\u00d73=492\", \"758\u00d78=6064\"),\n    @(\"557\u00d74=2228\", \"587\u00d75=2935\"),\n    @(\"628\u00d73=1884\", \"723\u00d72=1446\"),\n    @(\"697\u00d79=6273\", \"611\u00d72=1222\"),\n    @(\"433\u00d76=2598\", \"376\u00d72=752\"),\n    @(\"207\u00d79=1863\", \"774\u00d72=1548\"),\n    @(\"485\u00d73=1455\", \"372\u00d73=1116\"),\n    @(\"519\u00d72=1038\", \"382\u00d76=2292\"),\n    @(\"982\u00d78=7856\", \"906\u00d77=6342\"),\n    @(\"298\u00d73=894\", \"810\u00d74=3240\"),\n    @(\"369\u00d76=2214\", \"298\u00d76=1788\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
